$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 86 (data rows shift down by 2).
$ws.Range("A86:R87").EntireRow.Insert()

# Fill in the two newly inserted rows with their data.
# New row 86
$ws.Cells.Item(86, 1).Value = 6
$ws.Cells.Item(86, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(86, 3).Value = "Metropolitana"
$ws.Cells.Item(86, 4).Value = 44518
$ws.Cells.Item(86, 5).Value = 13
$ws.Cells.Item(86, 6).Value = 100112001
$ws.Cells.Item(86, 7).Value = "Berenjena"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 250
$ws.Cells.Item(86, 11).Value = 7000
$ws.Cells.Item(86, 12).Value = 8000
$ws.Cells.Item(86, 13).Value = 7480
$ws.Cells.Item(86, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(86, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value = 150
$ws.Cells.Item(86, 17).Value = 50
$ws.Cells.Item(86, 18).Value = "Hortaliza"

# New row 87
$ws.Cells.Item(87, 1).Value = 6
$ws.Cells.Item(87, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(87, 3).Value = "Metropolitana"
$ws.Cells.Item(87, 4).Value = 44518
$ws.Cells.Item(87, 5).Value = 13
$ws.Cells.Item(87, 6).Value = 100112001
$ws.Cells.Item(87, 7).Value = "Berenjena"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 150
$ws.Cells.Item(87, 11).Value = 12000
$ws.Cells.Item(87, 12).Value = 13000
$ws.Cells.Item(87, 13).Value = 12467
$ws.Cells.Item(87, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(87, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(87, 16).Value = 208
$ws.Cells.Item(87, 17).Value = 60
$ws.Cells.Item(87, 18).Value = "Hortaliza"
